# Updated cryptos list on Wed Jul 19 20:37:32 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures scraped
# from coinranking.com, and fixes the Quant/Aptos row ordering (rows 46-47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as plain text in the sheet (e.g. "7.000" or
# "30.047.32"), so a leading text-qualifier apostrophe is used whenever the
# new value would otherwise be auto-recognized by Excel as a number, which
# would silently drop formatting such as trailing zeros.

$ws.Range("D2").Value = '30.047.32'
$ws.Range("E2").Value = '  +0.93%  '

$ws.Range("D3").Value = '1.905.95'
$ws.Range("E3").Value = '  +0.55%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''0.8309'
$ws.Range("E5").Value = '  +9.23%  '

$ws.Range("D6").Value = '''242.05'
$ws.Range("E6").Value = '  +0.85%  '

$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '''0.3237'
$ws.Range("E8").Value = '  +6.75%  '

$ws.Range("D9").Value = '''26.79'
$ws.Range("E9").Value = '  +5.90%  '

$ws.Range("D10").Value = '''0.07035'
$ws.Range("E10").Value = '  +3.47%  '

$ws.Range("D11").Value = '''0.08034'
$ws.Range("E11").Value = '  +0.81%  '

$ws.Range("D12").Value = '''0.7514'
$ws.Range("E12").Value = '  +2.36%  '

$ws.Range("D13").Value = '1.910.63'
$ws.Range("E13").Value = '  +1.13%  '

$ws.Range("D14").Value = '''5.220'
$ws.Range("E14").Value = '  +1.47%  '

$ws.Range("D15").Value = '''92.84'
$ws.Range("E15").Value = '  +2.34%  '

$ws.Range("D16").Value = '30.061.73'
$ws.Range("E16").Value = '  +0.98%  '

$ws.Range("D17").Value = '''14.16'
$ws.Range("E17").Value = '  +2.81%  '

$ws.Range("D18").Value = '''5.946'
$ws.Range("E18").Value = '  +0.91%  '

$ws.Range("D19").Value = '''244.51'
$ws.Range("E19").Value = '  +1.53%  '

$ws.Range("E20").Value = '  +1.30%  '

$ws.Range("D21").Value = '2.161.09'
$ws.Range("E21").Value = '  +1.03%  '

$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").Value = '''1.001'
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").Value = '''7.000'
$ws.Range("E24").Value = '  +1.64%  '

$ws.Range("E25").Value = '  +23.86%  '

$ws.Range("D26").Value = '''168.66'
$ws.Range("E26").Value = '  +1.34%  '

$ws.Range("D27").Value = '''9.242'
$ws.Range("E27").Value = '  +0.47%  '

$ws.Range("D28").Value = '''18.98'

$ws.Range("D29").Value = '''2.092'
$ws.Range("E29").Value = '  +3.66%  '

$ws.Range("E30").Value = '  -1.73%  '

$ws.Range("D31").Value = '''1.517'
$ws.Range("E31").Value = '  +0.27%  '

$ws.Range("D32").Value = '''4.287'

$ws.Range("D33").Value = '''0.05713'
$ws.Range("E33").Value = '  +10.14%  '

$ws.Range("D34").Value = '''4.092'
$ws.Range("E34").Value = '  +0.89%  '

$ws.Range("E35").Value = '  +3.75%  '

$ws.Range("D36").Value = '''0.7343'
$ws.Range("E36").Value = '  +1.65%  '

$ws.Range("D37").Value = '''2.731'
$ws.Range("E37").Value = '  +0.65%  '

$ws.Range("D38").Value = '''0.01909'
$ws.Range("E38").Value = '  -0.12%  '

$ws.Range("D39").Value = '''2.798'
$ws.Range("E39").Value = '  +0.98%  '

$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("D41").Value = '''72.36'
$ws.Range("E41").Value = '  +0.78%  '

$ws.Range("D42").Value = '''5.978'
$ws.Range("E42").Value = '  -2.51%  '

$ws.Range("D43").Value = '''0.8431'
$ws.Range("E43").Value = '  +1.95%  '

$ws.Range("D44").Value = '''1.001'
$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("D45").Value = '''1.901'
$ws.Range("E45").Value = '  +1.36%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '''7.622'
$ws.Range("E46").Value = '  +0.44%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''101.30'
$ws.Range("E47").Value = '  +1.84%  '

$ws.Range("D48").Value = '''9.731'
$ws.Range("E48").Value = '  +0.44%  '

$ws.Range("D49").Value = '''994.37'
$ws.Range("E49").Value = '  +9.56%  '

$ws.Range("D50").Value = '2.069.49'
$ws.Range("E50").Value = '  +1.39%  '

$ws.Range("D51").Value = '''36.36'
$ws.Range("E51").Value = '  +1.03%  '
